$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Department column (C) values changed from "FACULTY OF ENGLISH" to "English"
$ws.Range("C2").Value = "English"
$ws.Range("C3").Value = "English"

# S3 previously duplicated the promotionValidity text from R3; it is now cleared
$ws.Range("S3").Value = ""
